$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(699640,699827,701350,701747,702499,702850,702951,703205,703347,703400,703495,703497,703498,703740,703759)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("A1").Select()
